# Daily attendance processing - 2025-12-19 12:46:00
#
# The "Recorded By" column (G) lists the accounts that touched each
# attendance record as a comma-separated string, e.g.
#   "dnasr281@gmail.com, System"
# Re-order each list so that any "System" / "system" entries come first,
# followed by the remaining (human/service) accounts, preserving their
# relative order. Cells that don't contain a "System" token are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $orig = $cell.Text

    if ($orig -eq "") { continue }

    $parts = $orig -split ", "

    $sysParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p -eq "System" -or $p -eq "system") {
            $sysParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($sysParts.Count -eq 0) { continue }

    $newParts = $sysParts + $otherParts
    $newVal = $newParts -join ", "

    if ($newVal -ne $orig) {
        $cell.Value = $newVal
    }
}
